$d = $word.ActiveDocument

# The ESUP permit template's signature block used to show a static label
# ("Inspector of Mines") under the signature line. Replace it with the
# Carbone merge-field placeholder so the generated permit prints the
# actual issuing inspector's name: {d.issuing_inspector_name}
#
# Build the replacement the same way the source template spells out its
# other {d....} fields elsewhere in the document -- as the three pieces
# "{d." / "issuing_inspector_name" / "}" -- while keeping the run's
# existing character formatting (Arial complex-script font, -10 spacing,
# 12 half-point monospaced "fine print" size) untouched.
$rng = $d.Content
$rng.Find.Execute("Inspector of Mines", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)

if ($rng.Find.Found) {
    $rng.Text = "{d."
    $rng.Collapse(0)
    $rng.InsertAfter("issuing_inspector_name")
    $rng.Collapse(0)
    $rng.InsertAfter("}")
}
